# Apply "contingencies with rene fine" edit to lines_states sheet.
# Two new line rows (line7, line8) are inserted after line6, so the
# existing extr1..extr8 rows shift down by two rows (8-15 -> 10-17).
# New contingency data (columns C/D/E) replaces the old values for every
# row from line7 down through extr8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the extr1..extr8 block (old rows 8-15) down to rows 10-17,
# working bottom-to-top so we never clobber a row before reading it.
# Value2 is used (rather than Value) so text cells copy as plain text,
# not a stringified COM property descriptor.
for ($r = 15; $r -ge 8; $r--) {
    $dest = $r + 2
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($dest, $c).Value2 = $ws.Cells.Item($r, $c).Value2
    }
}

# Rows 16 and 17 are brand-new cells beyond the original A1:E15 range,
# so column A needs its bordered/bold/centered style applied explicitly
# (copy format only, so no stray unused style gets minted).
$ws.Range("A7").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

# Full target table (row -> A, name, C, D, E) after the shift.
$data = @(
    @{ Row = 8;  A = 6;  Name = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  Name = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  Name = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  Name = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; Name = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; Name = "extr4"; C = 7;  D = 8;  E = $true  },
    @{ Row = 14; A = 12; Name = "extr5"; C = 9;  D = 11; E = $true  },
    @{ Row = 15; A = 13; Name = "extr6"; C = 7;  D = 11; E = $false },
    @{ Row = 16; A = 14; Name = "extr7"; C = 5;  D = 7;  E = $false },
    @{ Row = 17; A = 15; Name = "extr8"; C = 8;  D = 5;  E = $false }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.Name
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
}
